$wb = $excel.ActiveWorkbook
$originalActiveSheet = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet.
#    "总计" is currently the 4th sheet. We duplicate the "2021-Q4" sheet
#    (same A:H layout/formatting used by every per-quarter sheet) and drop
#    it in front of "总计" -- this keeps sheetPr/outline settings and all
#    the header/column styling intact, then we overwrite its data.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(4)
$srcSheet = $wb.Worksheets.Item(3)
$srcSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item(4)
$newSheet.Name = "2022-Q1"

# The source sheet (2021-Q4) has 34 data rows; the new sheet only needs 6
# (rows 2-7), so drop the now-unused tail (also shrinks "dimension" back
# down to A1:H7, matching the new row count).
$newSheet.Range("A8:H34").Clear()

# Header row (overwrite the copied "2021-Q4" header text)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row index column
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4
$newSheet.Range("A7").Value = 5

# Fund holdings data.  B..G are stored as plain text (matching the rest of
# the workbook's per-quarter sheets) -- a leading apostrophe forces
# text-entry so things like "001766" keep their leading zero and "10.35"
# is not coerced into a numeric cell. The style is reset to "Normal"
# afterwards so no stray number-format style lingers on the cell.
$newSheet.Range("B2").Value = "'001766"
$newSheet.Range("C2").Value = "上投摩根医疗健康股票"
$newSheet.Range("D2").Value = "'10.35"
$newSheet.Range("E2").Value = "'80.54"
$newSheet.Range("F2").Value = "'2.67"
$newSheet.Range("G2").Value = "'0.2763"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "'009468"
$newSheet.Range("C3").Value = "博时健康成长主题双周定期可赎回混合A"
$newSheet.Range("D3").Value = "'6.39"
$newSheet.Range("E3").Value = "'85.57"
$newSheet.Range("F3").Value = "'3.83"
$newSheet.Range("G3").Value = "'0.2447"
$newSheet.Range("H3").Value = 4

$newSheet.Range("B4").Value = "'009469"
$newSheet.Range("C4").Value = "博时健康成长主题双周定期可赎回混合C"
$newSheet.Range("D4").Value = "'1.11"
$newSheet.Range("E4").Value = "'85.57"
$newSheet.Range("F4").Value = "'3.83"
$newSheet.Range("G4").Value = "'0.0425"
$newSheet.Range("H4").Value = 4

$newSheet.Range("B5").Value = "'011214"
$newSheet.Range("C5").Value = "招商惠润一年定期开放混合型发起式管理人中管理人（MOM）证券投资基金A"
$newSheet.Range("D5").Value = "'0.67"
$newSheet.Range("E5").Value = "'81.20"
$newSheet.Range("F5").Value = "'2.33"
$newSheet.Range("G5").Value = "'0.0156"
$newSheet.Range("H5").Value = 6

$newSheet.Range("B6").Value = "'004536"
$newSheet.Range("C6").Value = "嘉实中小企业量化活力灵活配置混合"
$newSheet.Range("D6").Value = "'0.17"
$newSheet.Range("E6").Value = "'90.06"
$newSheet.Range("F6").Value = "'2.58"
$newSheet.Range("G6").Value = "'0.0044"
$newSheet.Range("H6").Value = 1

$newSheet.Range("B7").Value = "'011215"
$newSheet.Range("C7").Value = "招商惠润一年定期开放混合型发起式管理人中管理人（MOM）证券投资基金C"
$newSheet.Range("D7").Value = "'0.09"
$newSheet.Range("E7").Value = "'81.20"
$newSheet.Range("F7").Value = "'2.33"
$newSheet.Range("G7").Value = "'0.0021"
$newSheet.Range("H7").Value = 6

# Column B (fund code) and D:G were forced to text via the leading
# apostrophe above; clear the resulting quote-prefix style so the cells
# fall back to the sheet's default (unstyled) look, same as column C.
$newSheet.Range("B2:B7").Style = "Normal"
$newSheet.Range("D2:G7").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a 2022-Q1 summary row above the existing
#    quarters (shifting 2021-Q4/Q3/Q2 down by one row).
# ---------------------------------------------------------------------------
$ws_total = $wb.Worksheets.Item(5)

# Make sure row 5 (newly used) carries the same column-A style (s="2") as the
# rows above it before we shift values into it.
$ws_total.Range("A2:A4").Copy()
$ws_total.Range("A3:A5").PasteSpecial(-4122)

for ($r = 4; $r -ge 2; $r--) {
  $newR = $r + 1
  $ws_total.Cells.Item($newR, 2).Value = $ws_total.Cells.Item($r, 2).Value()
  $ws_total.Cells.Item($newR, 3).Value = $ws_total.Cells.Item($r, 3).Value()
  $ws_total.Cells.Item($newR, 4).Value = $ws_total.Cells.Item($r, 4).Value()
}

# Column A is just the running 0-based row index -- renumber it in place.
$ws_total.Range("A2").Value = 0
$ws_total.Range("A3").Value = 1
$ws_total.Range("A4").Value = 2
$ws_total.Range("A5").Value = 3

$ws_total.Range("B2").Value = "2022-Q1"
$ws_total.Range("C2").Value = 6
$ws_total.Range("D2").Value = 0.59

# Restore whichever sheet/tab was active before this script ran (adding /
# copying into sheets shifts Excel's "active sheet" as a side effect).
$originalActiveSheet.Activate()
